$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-5 and 7 per recomputed s_vals
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 3
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 1
$ws.Range("G7").Value = 1
